# Auto-generated Excel COM-interop script to apply odds updates
# from the Jogos_da_Semana_FlashScore_2024-10-24 diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 1.79
    "G3" = 1.41
    "O3" = 1.36
    "P3" = 3
    "Q3" = 2.1
    "R3" = 1.7
    "G4" = 1.95
    "H4" = 3.3
    "I4" = 4
    "J4" = 2.63
    "L4" = 4.75
    "M4" = 1.08
    "N4" = 8
    "Q4" = 2.3
    "R4" = 1.6
    "U4" = 2
    "V4" = 1.69
    "X4" = 8.5
    "Z4" = 17
    "AE4" = 17
    "AH4" = 9.5
    "AI4" = 19
    "AK4" = 41
    "AL4" = 34
    "AO4" = 11
    "AQ4" = 41
    "AX4" = 23
    "G5" = 1.17
    "H5" = 6.25
    "I5" = 19
    "J5" = 1.62
    "L5" = 13
    "Q5" = 1.8
    "R5" = 2
    "U5" = 3
    "V5" = 1.33
    "W5" = 5.5
    "Z5" = 6
    "AC5" = 11
    "AD5" = 13
    "AE5" = 41
    "AF5" = 151
    "AH5" = 29
    "AI5" = 81
    "AK5" = 351
    "AL5" = 151
    "AM5" = 151
    "AN5" = 2.88
    "AO5" = 5
    "AP5" = 23
    "AS5" = 251
    "AU5" = 13
    "AX5" = 67
    "AY5" = 67
    "V6" = 1.63
    "G7" = 2.4
    "H7" = 3.3
    "I7" = 3
    "J7" = 3.1
    "M7" = 1.07
    "N7" = 8.5
    "Q7" = 2.15
    "R7" = 1.63
    "S7" = 1.44
    "T7" = 2.63
    "U7" = 1.87
    "V7" = 1.87
    "Y7" = 9.5
    "AC7" = 8.5
    "AD7" = 6
    "AL7" = 26
    "AN7" = 4.33
    "AO7" = 13
    "AP7" = 26
    "AY7" = 29
    "L8" = 2.62
    "S8" = 1.33
    "G9" = 3.5
    "H9" = 3.6
    "I9" = 1.9
    "K9" = 2.25
    "L9" = 2.4
    "Q9" = 1.57
    "R9" = 2.12
    "S9" = 1.3
    "T9" = 3.32
    "U9" = 1.52
    "V9" = 2.22
    "W9" = 13.5
    "X9" = 22
    "Z9" = 50
    "AC9" = 14
    "AD9" = 7.3
    "AF9" = 40
    "AH9" = 9.75
    "AI9" = 11
    "AK9" = 18
    "AL9" = 13.5
    "AM9" = 19.5
    "AP9" = 22
    "AT9" = 3.15
    "AW9" = 4.05
    "AX9" = 9.25
    "AY9" = 15
    "AZ9" = 30
    "BA9" = 50
    "BB9" = 150
    "M10" = 1.06
    "N10" = 10
    "Q10" = 1.97
    "R10" = 1.77
    "S10" = 1.4
    "S11" = 1.44
    "T11" = 2.63
    "G13" = 2.2
    "I13" = 3.2
    "J13" = 2.75
    "L13" = 3.6
    "M13" = 1.03
    "N13" = 10
    "AA13" = 17
    "AH13" = 11
    "AI13" = 17
    "AK13" = 34
    "AL13" = 26
    "M14" = 1.03
    "N14" = 15
    "S16" = 1.5
    "T16" = 2.37
    "S17" = 1.54
    "S18" = 1.58
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

